$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all the old vendor/appliance/brand/region/etc. columns (C:K), keeping only
# what will become the new A (Pincode) / B (Appliance) layout
$ws.Columns("C:K").Delete()

# New header row: Pincode, Appliance
$ws.Range("A1").Value = "Pincode"
$ws.Range("B1").Value = "Appliance"

# New data row: the merge-field placeholders
$ws.Range("A2").Value = "{vendor:Pincode}"
$ws.Range("B2").Value = "{vendor:Appliance}"

# Center-align both rows (header keeps its bold font, data row gets the default font)
$ws.Range("A1:B1").HorizontalAlignment = -4108
$ws.Range("A2:B2").HorizontalAlignment = -4108

# Resize the two remaining columns
$ws.Columns.Item(1).ColumnWidth = 27.333333333333332
$ws.Columns.Item(2).ColumnWidth = 98.83333333333333

# Reset the selection to the top-left cell
$ws.Range("A1").Select()
